# Setup OrangeHRM locally & Added DB Verification
#
# - Insert a new "empVerification" sheet between "validLoginData" and
#   "inValidLoginData" with employee-id/name verification data.
# - Update the valid-login credentials on "validLoginData".
# - Update the invalid-login password on "inValidLoginData".

$wb = $excel.ActiveWorkbook

$validSheet = $wb.Worksheets.Item("validLoginData")
$invalidSheet = $wb.Worksheets.Item("inValidLoginData")

# --- Update validLoginData credentials (row 2) ---
$validSheet.Range("A2").Value = "orangehrm_Subhasis"
$validSheet.Range("B2").Value = "Spal@HRM7"

# --- Update inValidLoginData password (row 2) ---
$invalidSheet.Range("B2").Value = "admin123"
[void]$invalidSheet.Range("A2:B2").Select()

# --- Insert new "empVerification" sheet right before inValidLoginData ---
# (Copying an existing sheet keeps the worksheet's namespace/markup-
# compatibility declarations consistent with its siblings; Copy() itself
# does not return the new sheet, and the old $invalidSheet reference keeps
# its stale pre-copy Index, so re-fetch it live by name to find the copy.)
[void]$invalidSheet.Copy($invalidSheet, $null)
$invalidSheet = $wb.Worksheets.Item("inValidLoginData")
$empSheet = $wb.Worksheets.Item($invalidSheet.Index - 1)
$empSheet.Name = "empVerification"

$empSheet.Range("A1").Value = "empl_id"
$empSheet.Range("B1").Value = "emp_name"
$empSheet.Range("A2").Value = 1
$empSheet.Range("B2").Value = "Subh"
[void]$empSheet.Range("B2").Select()
